$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: a leading apostrophe in a cell value is interpreted by Excel as a
# "force text / quote-prefix" marker and is stripped from the stored text
# (while also flipping on the cell's quotePrefix style). To get a literal
# leading apostrophe into the stored string - and to leave the cell's
# style untouched - we double the leading apostrophe (only the first is
# consumed as the marker) and then explicitly reset the style back to
# Normal afterwards.

$data = @(
  @{ Row = 2;  B = "''Bacteroides_cellulosilyticus_DSM_14838.mat'"; C = 0 },
  @{ Row = 3;  B = "''Bacteroides_coprocola_M16_DSM_17136.mat'"; C = 0 },
  @{ Row = 4;  B = "''Bacteroides_fluxus_YIT_12057.mat'"; C = 0 },
  @{ Row = 5;  B = "''Bacteroides_oleiciplenus_YIT_12058.mat'"; C = 0 },
  @{ Row = 6;  B = "''Bacteroides_ovatus_ATCC_8483.mat'"; C = 0 },
  @{ Row = 7;  B = "''Bacteroides_salyersiae_WAL_10018.mat'"; C = 0 },
  @{ Row = 8;  B = "''Bacteroides_stercoris_ATCC_43183.mat'"; C = 0 },
  @{ Row = 9;  B = "''Bacteroides_thetaiotaomicron_VPI_5482.mat'"; C = 0.017 },
  @{ Row = 10; B = "''Bacteroides_uniformis_ATCC_8492.mat'"; C = 0 },
  @{ Row = 11; B = "''Bacteroides_vulgatus_ATCC_8482.mat'"; C = 0.011 },
  @{ Row = 12; B = "''Bifidobacterium_animalis_lactis_AD011.mat'"; C = 0 },
  @{ Row = 13; B = "''Enterococcus_faecalis_OG1RF_ATCC_47077.mat'"; C = 0 },
  @{ Row = 14; B = "''Flavonifractor_plautii_ATCC_29863.mat'"; C = 0 },
  @{ Row = 15; B = "''Gordonibacter_pamelaeae_7_10_1_bT_DSM_19378.mat'"; C = 0.011 },
  @{ Row = 16; B = "''Odoribacter_laneus_YIT_12061.mat'"; C = 0.147 },
  @{ Row = 17; B = "''Parabacteroides_distasonis_ATCC_8503.mat'"; C = 0 },
  @{ Row = 18; B = "''Parabacteroides_johnsonii_DSM_18315.mat'"; C = 0.8139999999999999 }
)

foreach ($item in $data) {
  $cellB = $ws.Cells.Item($item.Row, 2)
  $cellB.Value = $item.B
  $cellB.Style = "Normal"
  $ws.Cells.Item($item.Row, 3).Value = $item.C
}
